$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 285
$ws.Range("I31").Value = 285
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 855
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -625

$ws.Range("H33").Value = 299.18182
$ws.Range("I33").Value = 269.05264
$ws.Range("J33").Value = 490
$ws.Range("K33").Value = 269.05264
$ws.Range("L33").Value = 490
$ws.Range("M33").Value = -40.05264

$ws.Range("H76").Value = 3550
$ws.Range("I76").Value = 3550
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3550
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3235
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 3550
$ws.Range("I79").Value = 3550
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3550
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2458
$ws.Range("N79").ClearContents()

$ws.Range("H98").Value = 722.0952
$ws.Range("I98").Value = 704.7778
$ws.Range("J98").Value = 826
$ws.Range("K98").Value = 704.7778
$ws.Range("L98").Value = 826
$ws.Range("M98").Value = 793.2222
$ws.Range("N98").Value = -3822

$ws.Range("H104").Value = 621.8
$ws.Range("I104").Value = 702.25
$ws.Range("J104").Value = 300
$ws.Range("K104").Value = 2106.75
$ws.Range("L104").Value = 900
$ws.Range("M104").Value = -359.75

$ws.Range("H107").Value = 1067.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1067.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1067.5
$ws.Range("N107").Value = -4907.5

$ws.Range("H113").Value = 3663.28
$ws.Range("I113").Value = 3511.4614
$ws.Range("J113").Value = 3827.75
$ws.Range("K113").Value = 3511.4614
$ws.Range("L113").Value = 3827.75
$ws.Range("M113").Value = -257.4614000000001
$ws.Range("N113").Value = -10335.75

$ws.Range("H119").Value = 2300
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 2300
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 6900
$ws.Range("N119").Value = -16576

$ws.Range("H122").Value = 722.0952
$ws.Range("I122").Value = 704.7778
$ws.Range("J122").Value = 826
$ws.Range("K122").Value = 2114.3334
$ws.Range("L122").Value = 2478
$ws.Range("M122").Value = 335.6666
$ws.Range("N122").Value = -7378

$ws.Range("H125").Value = 894
$ws.Range("I125").Value = 265
$ws.Range("J125").Value = 1313.3334
$ws.Range("K125").Value = 2385
$ws.Range("L125").Value = 11820.0006
$ws.Range("M125").Value = 75
$ws.Range("N125").Value = -16740.0006

$ws.Range("H132").Value = 2708.9744
$ws.Range("I132").Value = 2495.238
$ws.Range("J132").Value = 3606.6667
$ws.Range("K132").Value = 7485.714
$ws.Range("L132").Value = 10820.0001
$ws.Range("M132").Value = -4955.714
$ws.Range("N132").Value = -15880.0001

$ws.Range("H138").Value = 2669038.8
$ws.Range("I138").Value = 1486.7742
$ws.Range("J138").Value = 4548450.5
$ws.Range("K138").Value = 4460.3226
$ws.Range("L138").Value = 13645351.5
$ws.Range("M138").Value = 679.6773999999996
$ws.Range("N138").Value = -13655631.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5643.29
$ws.Range("I32").Value = 5244.125
$ws.Range("J32").Value = 8570.5
$ws.Range("K32").Value = 5244.125
$ws.Range("L32").Value = 8570.5
$ws.Range("M32").Value = -4957.125
$ws.Range("N32").Value = -9144.5

$ws.Range("H132").Value = 9475272
$ws.Range("I132").Value = 11652685
$ws.Range("J132").Value = 112397.8
$ws.Range("K132").Value = 34958055
$ws.Range("L132").Value = 337193.4
$ws.Range("M132").Value = -34955525

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2058.1428
$ws.Range("I20").Value = 1793.3334
$ws.Range("J20").Value = 2256.75
$ws.Range("K20").Value = 1793.3334
$ws.Range("L20").Value = 2256.75
$ws.Range("M20").Value = -1546.3334
$ws.Range("N20").Value = -2750.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 26317806
$ws.Range("I58").Value = 38463284
$ws.Range("J58").Value = 2598.75
$ws.Range("K58").Value = 38463284
$ws.Range("L58").Value = 2598.75
$ws.Range("M58").Value = -38463081

$ws.Range("H136").Value = 26317806
$ws.Range("I136").Value = 38463284
$ws.Range("J136").Value = 2598.75
$ws.Range("K136").Value = 115389852
$ws.Range("L136").Value = 7796.25
$ws.Range("M136").Value = -115387302

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 278038.84
$ws.Range("I11").Value = 189.7
$ws.Range("J11").Value = 625350.25
$ws.Range("K11").Value = 569.0999999999999
$ws.Range("L11").Value = 1876050.75
$ws.Range("M11").Value = -429.0999999999999
$ws.Range("N11").Value = -1876330.75

$ws.Range("H12").Value = 102.457146
$ws.Range("I12").Value = 60.75
$ws.Range("J12").Value = 137.57895
$ws.Range("K12").Value = 182.25
$ws.Range("L12").Value = 412.73685
$ws.Range("M12").Value = -9.25
$ws.Range("N12").Value = -758.73685

$ws.Range("H54").Value = 2002.5
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 2002.5
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 6007.5
$ws.Range("N54").Value = -7125.5

$ws.Range("H109").Value = 3175
$ws.Range("I109").Value = 1100
$ws.Range("J109").Value = 3866.6667
$ws.Range("K109").Value = 3300
$ws.Range("L109").Value = 11600.0001
$ws.Range("M109").Value = -2260
$ws.Range("N109").Value = -13680.0001

$ws.Range("H126").Value = 3168.25
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3168.25
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 9504.75
$ws.Range("N126").Value = -19384.75
$ws.Range("M126").ClearContents()

$ws.Range("H141").Value = 14103.333
$ws.Range("I141").Value = 6832.5
$ws.Range("J141").Value = 19920
$ws.Range("K141").Value = 20497.5
$ws.Range("L141").Value = 59760
$ws.Range("M141").Value = -15317.5
$ws.Range("N141").Value = -70120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H70").Value = 36293.72
$ws.Range("I70").Value = 52480.953
$ws.Range("J70").Value = 5390.8184
$ws.Range("K70").Value = 52480.953
$ws.Range("L70").Value = 5390.8184
$ws.Range("M70").Value = -52210.953

$ws.Range("H73").Value = 36293.72
$ws.Range("I73").Value = 52480.953
$ws.Range("J73").Value = 5390.8184
$ws.Range("K73").Value = 52480.953
$ws.Range("L73").Value = 5390.8184
$ws.Range("M73").Value = -51544.953

$ws.Range("H80").Value = 3984.0588
$ws.Range("I80").Value = 3001.25
$ws.Range("J80").Value = 4286.4614
$ws.Range("K80").Value = 3001.25
$ws.Range("L80").Value = 4286.4614
$ws.Range("M80").Value = -2003.25
$ws.Range("N80").Value = -6282.4614

$ws.Range("H83").Value = 3984.0588
$ws.Range("I83").Value = 3001.25
$ws.Range("J83").Value = 4286.4614
$ws.Range("K83").Value = 15006.25
$ws.Range("L83").Value = 21432.307
$ws.Range("M83").Value = -10014.25
$ws.Range("N83").Value = -31416.307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H93").Value = 1149.7333
$ws.Range("I93").Value = 1138.7826
$ws.Range("J93").Value = 1185.7142
$ws.Range("K93").Value = 1138.7826
$ws.Range("L93").Value = 1185.7142
$ws.Range("M93").Value = 109.2174

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 10002006
$ws.Range("I20").Value = 20000010
$ws.Range("J20").Value = 4001
$ws.Range("K20").Value = 20000010
$ws.Range("L20").Value = 4001
$ws.Range("M20").Value = -19999770

$ws.Range("H92").Value = 40549.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 40549.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 40549.5
$ws.Range("N92").Value = -45541.5

$ws.Range("H132").Value = 36006.035
$ws.Range("I132").Value = 34385.066
$ws.Range("J132").Value = 37742.785
$ws.Range("K132").Value = 103155.198
$ws.Range("L132").Value = 113228.355
$ws.Range("M132").Value = -100625.198
$ws.Range("N132").Value = -118288.355
